$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B2 formula to use a relative path instead of the hard-coded
# absolute user path (the redundant/broken path that triggered the
# "browser update issue").
$ws.Range("B2").Formula = '="samples\"&A2&" "&C2'

# Move the active selection from B4 to B2.
$ws.Range("B2").Select()
